# Weekly data refresh: a new week of price data is inserted at row 37
# (the start of this product/market block), pushing every subsequent
# row's varying fields down by one row. The former last row (114) is
# preserved by appending it as the new row 115.
#
# Columns that stay constant for every data row in this sheet
# (A,B,C,E,F,G,H,I,N,Q,R) are left untouched for existing rows; only the
# columns that actually vary row-to-row (D,J,K,L,M,O,P) are shifted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$varyCols = @("D", "J", "K", "L", "M", "O", "P")

# Row 115 is brand new, so the columns that are constant down the whole
# block (A,B,C,E,F,G,H,I,N,Q,R) need to be populated explicitly there -
# every other destination row already carries them from its own prior
# data.
$constCols = @("A", "B", "C", "E", "F", "G", "H", "I", "N", "Q", "R")
foreach ($col in $constCols) {
    $ws.Range($col + "115").Value = $ws.Range($col + "114").Value2()
}

# Shift rows 37..114 down to rows 38..115, working from the bottom up so
# we never overwrite a source row before it has been copied. Value2 is
# used (instead of Value) so a date-formatted cell's raw serial number is
# copied verbatim instead of Excel re-boxing it as a Date and silently
# reformatting the brand-new row 115 cells.
for ($r = 115; $r -ge 38; $r--) {
    $src = $r - 1
    foreach ($col in $varyCols) {
        $ws.Range($col + "$r").Value = $ws.Range($col + "$src").Value2()
    }
}

# Row 115's date cell (D) had no prior style of its own (it didn't exist
# before), so give it the same date format used by every other cell in
# column D.
$ws.Range("D115").NumberFormat = $ws.Range("D114").NumberFormat()

# Row 37 becomes the newly reported week's data.
$ws.Range("D37").Value = 45162
$ws.Range("J37").Value = 16
$ws.Range("K37").Value = 21000
$ws.Range("L37").Value = 21000
$ws.Range("M37").Value = 21000
$ws.Range("O37").Value = "Región Metropolitana"
$ws.Range("P37").Value = 7000
